{"js": "// Update the single-column results table: the first 12 rows get refreshed\n// summary-statistic values, and the last 3 rows (which held full per-run\n// tab-separated data dumps) get collapsed down to a single summary value\n// each, matching the new README/docx stats preparation.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// row index (0-based) -> new cell text\nconst newValues = {\n  0: \"0M\",\n  1: \"0M\",\n  2: \"0M\",\n  3: \"319\",\n  4: \"0.00001\",\n  5: \"0.00064\",\n  6: \"0.00017\",\n  7: \"0.00006\",\n  8: \"0.00032\",\n  9: \"0.00037\",\n  10: \"0.00040\",\n  11: \"0.06493\",\n  43: \"99.67\",\n  44: \"0.06\",\n  45: \"19\",\n};\n\nfor (const rowIndex of Object.keys(newValues)) {\n  const cell = table.getCellOrNullObject(Number(rowIndex), 0);\n  cell.value = newValues[rowIndex];\n}\n\nawait context.sync();\n", "ps1": "# Update the single-column results table: the first 12 rows get refreshed\n# summary-statistic values, and the last 3 rows (which held full per-run\n# tab-separated data dumps) get collapsed down to a single summary value\n# each, matching the new README/docx stats preparation.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$t.Cell(1, 1).Range.Text = \"0M\"\n$t.Cell(2, 1).Range.Text = \"0M\"\n$t.Cell(3, 1).Range.Text = \"0M\"\n$t.Cell(4, 1).Range.Text = \"319\"\n$t.Cell(5, 1).Range.Text = \"0.00001\"\n$t.Cell(6, 1).Range.Text = \"0.00064\"\n$t.Cell(7, 1).Range.Text = \"0.00017\"\n$t.Cell(8, 1).Range.Text = \"0.00006\"\n$t.Cell(9, 1).Range.Text = \"0.00032\"\n$t.Cell(10, 1).Range.Text = \"0.00037\"\n$t.Cell(11, 1).Range.Text = \"0.00040\"\n$t.Cell(12, 1).Range.Text = \"0.06493\"\n$t.Cell(44, 1).Range.Text = \"99.67\"\n$t.Cell(45, 1).Range.Text = \"0.06\"\n$t.Cell(46, 1).Range.Text = \"19\"\n"}
